$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (EEE 2113, Sec A): update course title and class times
$ws.Range("B3").Value = "Electrical Circuits AA"
$ws.Range("H3").Value = "06:00:PM - 06:50:PM"
$ws.Range("I3").Value = "06:00:PM - 06:50:PM"

# Row 5 (EEE 2123, Sec F): update course title (with "BB" in bold) and class times
$ws.Range("B5").Value = "Electronics BB"
$ws.Range("B5").Characters(1, 12).Font.Name = "Times New Roman"
$ws.Range("B5").Characters(1, 12).Font.Size = 12
$ws.Range("B5").Characters(13, 2).Font.Name = "Times New Roman"
$ws.Range("B5").Characters(13, 2).Font.Size = 12
$ws.Range("B5").Characters(13, 2).Font.Bold = $true
$ws.Range("H5").Value = "06:51:PM - 09:50:PM"
$ws.Range("I5").Value = "06:51:PM - 09:50:PM"
